# fix excel cell autosize
# Re-assign the tile letters in the "map" sheet (A1:N14). The underlying
# shared-string table swaps meaning of "B" (Black) and "W" (White) tiles,
# so every cell that used to resolve to "B" or "W" is rewritten here with
# its final, correct letter value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("map")

# Target grid, row by row (A..N), for rows 1..14
$grid = @(
    @("G","G","B","B","G","G","G","G","G","G","G","G","G","G"),
    @("G","G","G","G","B","B","G","G","B","B","G","G","G","G"),
    @("W","G","G","G","B","B","G","G","G","B","B","B","B","G"),
    @("W","G","G","B","B","B","B","G","B","B","B","B","G","G"),
    @("G","W","W","B","B","B","B","B","B","B","W","B","G","G"),
    @("B","B","B","W","B","B","B","B","B","B","G","G","W","G"),
    @("G","B","B","B","B","B","B","B","W","B","W","W","B","B"),
    @("B","G","B","B","B","B","B","B","B","B","B","B","B","G"),
    @("B","B","B","B","B","B","B","B","B","B","B","B","B","B"),
    @("G","G","G","B","B","B","W","B","B","B","W","G","B","G"),
    @("G","B","B","B","B","G","W","W","B","B","W","B","G","G"),
    @("B","G","B","B","B","G","G","G","W","B","G","G","G","B"),
    @("G","G","G","G","W","W","W","W","G","B","B","G","G","G"),
    @("G","G","B","B","W","W","W","G","G","G","G","G","G","G")
)

for ($r = 1; $r -le $grid.Length; $r++) {
    $rowValues = $grid[$r - 1]
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
